$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2025-11-27 Thursday" "2025-11-28 Friday"

Replace-Text "85×80=6800" "65×12=780"
Replace-Text "12×98=1176" "86×48=4128"
Replace-Text "33×75=2475" "64×60=3840"
Replace-Text "87×47=4089" "37×65=2405"
Replace-Text "88×14=1232" "57×40=2280"

Replace-Text "51×93=4743" "83×85=7055"
Replace-Text "49×70=3430" "33×73=2409"
Replace-Text "18×16=288" "14×80=1120"
Replace-Text "69×78=5382" "91×97=8827"
Replace-Text "43×70=3010" "36×62=2232"

Replace-Text "81×60=4860" "47×21=987"
Replace-Text "73×13=949" "27×92=2484"
Replace-Text "40×72=2880" "60×96=5760"
Replace-Text "68×11=748" "32×84=2688"
Replace-Text "88×20=1760" "26×63=1638"

Replace-Text "72×44=3168" "18×58=1044"
Replace-Text "45×28=1260" "63×87=5481"
Replace-Text "47×49=2303" "26×94=2444"
Replace-Text "44×64=2816" "64×97=6208"
Replace-Text "12×92=1104" "67×26=1742"

Replace-Text "32×98=3136" "79×39=3081"
Replace-Text "87×25=2175" "84×74=6216"
Replace-Text "72×22=1584" "38×65=2470"
Replace-Text "63×57=3591" "27×67=1809"
Replace-Text "17×56=952" "25×29=725"

Write-Output "Done"
